$d = $word.ActiveDocument

# 1) Remove the bullet-list paragraph "Ручное редактирование программ
#    тренировок" entirely, including its paragraph mark, so the remaining
#    paragraphs close the gap (the numbered list keeps only its first two
#    items).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "Ручное редактирование программ тренировок") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# 2) Tweak the closing sentence: "... могут быть изменены в информации об
#    аккаунте" becomes "... могут быть изменены в профиле".
$d.Content.Find.Execute("информации об аккаунте", $true, $false, $false, $false, $false, $true, 1, $false, "профиле", 2) | Out-Null
